$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Total Payment" label used to live in A9 (paired with a blank B9).
# Move/rename it to a bold "Total Amount" label in B9, and blank out A9.
$ws.Range("A9").Value() = " "
$ws.Range("B9").Value() = "Total Amount"
$ws.Range("B9").Font.Bold = $true

# Update the saved selection / active cell to B10:E10 (was B20).
[void]$ws.Range("B10:E10").Select()
